# Updated cryptos list on Wed Mar  6 16:07:33 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) holds plain-text values (e.g. "66.574.04"), not
# real numbers, so force text formatting before writing to avoid Excel
# auto-converting them to numeric literals; Style is reset straight after
# so no stray formatting is left behind on the cell.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "66.574.04"
$ws.Range("E2").Value = "  -0.11%  "

Set-TextValue $ws.Range("D3") "3.790.78"
$ws.Range("E3").Value = "  +0.37%  "

Set-TextValue $ws.Range("D4") "0.998"
$ws.Range("E4").Value = "  -1.05%  "

Set-TextValue $ws.Range("D5") "419.92"
$ws.Range("E5").Value = "  +0.91%  "

Set-TextValue $ws.Range("D6") "128.53"
$ws.Range("E6").Value = "  -7.21%  "

Set-TextValue $ws.Range("D7") "3.789.37"
$ws.Range("E7").Value = "  +0.52%  "

Set-TextValue $ws.Range("D8") "0.600"
$ws.Range("E8").Value = "  -5.26%  "

Set-TextValue $ws.Range("D9") "0.999"
$ws.Range("E9").Value = "  -0.02%  "

Set-TextValue $ws.Range("D10") "0.720"
$ws.Range("E10").Value = "  -3.79%  "

Set-TextValue $ws.Range("D11") "0.162"
$ws.Range("E11").Value = "  -3.58%  "

Set-TextValue $ws.Range("D12") "0.0000349"
$ws.Range("E12").Value = "  +11.02%  "

Set-TextValue $ws.Range("D13") "39.91"
$ws.Range("E13").Value = "  -7.66%  "

Set-TextValue $ws.Range("D14") "4.383.37"
$ws.Range("E14").Value = "  -0.63%  "

Set-TextValue $ws.Range("D15") "10.06"
$ws.Range("E15").Value = "  -0.83%  "

Set-TextValue $ws.Range("D16") "15.69"
$ws.Range("E16").Value = "  +18.15%  "

$ws.Range("E17").Value = "  -1.33%  "

Set-TextValue $ws.Range("D18") "3.788.82"
$ws.Range("E18").Value = "  +1.73%  "

Set-TextValue $ws.Range("D19") "19.40"
$ws.Range("E19").Value = "  -4.66%  "

Set-TextValue $ws.Range("D20") "66.646.77"
$ws.Range("E20").Value = "  -0.66%  "

Set-TextValue $ws.Range("D21") "1.07"
$ws.Range("E21").Value = "  -2.80%  "

Set-TextValue $ws.Range("D22") "402.32"
$ws.Range("E22").Value = "  -5.89%  "

Set-TextValue $ws.Range("D23") "14.21"
$ws.Range("E23").Value = "  -4.67%  "

Set-TextValue $ws.Range("D24") "83.30"
$ws.Range("E24").Value = "  -4.91%  "

Set-TextValue $ws.Range("D25") "2.98"
$ws.Range("E25").Value = "  -2.47%  "

Set-TextValue $ws.Range("D26") "36.91"
$ws.Range("E26").Value = "  +0.09%  "

$ws.Range("E27").Value = "  +10.58%  "

Set-TextValue $ws.Range("D28") "3.17"
$ws.Range("E28").Value = "  -3.23%  "

Set-TextValue $ws.Range("D29") "9.33"
$ws.Range("E29").Value = "  -2.57%  "

Set-TextValue $ws.Range("D30") "699.65"
$ws.Range("E30").Value = "  +0.22%  "

Set-TextValue $ws.Range("D31") "8.17"
$ws.Range("E31").Value = "  +16.13%  "

$ws.Range("E32").Value = "  +1.36%  "

Set-TextValue $ws.Range("D33") "12.26"
$ws.Range("E33").Value = "  -2.33%  "

Set-TextValue $ws.Range("D34") "0.120"
$ws.Range("E34").Value = "  -0.87%  "

$ws.Range("E35").Value = "  +0.21%  "

$ws.Range("E36").Value = "  -6.24%  "

$ws.Range("E37").Value = "  -8.40%  "

Set-TextValue $ws.Range("D38") "54.65"
$ws.Range("E38").Value = "  -3.98%  "

Set-TextValue $ws.Range("D39") "0.0₃0775"
$ws.Range("E39").Value = "  +23.82%  "

Set-TextValue $ws.Range("D40") "0.0450"
$ws.Range("E40").Value = "  -5.38%  "

$ws.Range("B41").Value = "NEARProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D41") "4.88"
$ws.Range("E41").Value = "  +13.58%  "

$ws.Range("B42").Value = "ThetaToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
Set-TextValue $ws.Range("D42") "2.93"
$ws.Range("E42").Value = "  -1.57%  "

$ws.Range("E43").Value = "  -0.06%  "

Set-TextValue $ws.Range("D44") "0.134"
$ws.Range("E44").Value = "  -5.84%  "

Set-TextValue $ws.Range("D45") "3.31"
$ws.Range("E45").Value = "  -2.40%  "

Set-TextValue $ws.Range("D46") "144.68"
$ws.Range("E46").Value = "  -1.27%  "

Set-TextValue $ws.Range("D47") "3.08"
$ws.Range("E47").Value = "  -1.86%  "

Set-TextValue $ws.Range("D48") "2.02"
$ws.Range("E48").Value = "  -3.39%  "

Set-TextValue $ws.Range("D49") "25.47"
$ws.Range("E49").Value = "  -1.97%  "

Set-TextValue $ws.Range("D50") "2.55"
$ws.Range("E50").Value = "  -0.59%  "

Set-TextValue $ws.Range("D51") "2.74"
$ws.Range("E51").Value = "  -2.98%  "
